$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.008", "11.80").
# Force text storage (matching the source workbook's inlineStr cells) so
# Excel's auto-type-detection on .Value does not coerce them to numbers and
# silently drop significant trailing zeros / multi-dot grouping.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D45","D46","D48","D49","D50","D51")
foreach ($ref in $dCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = '29.551.04'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.913.02'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.70%  '
$ws.Range("D5").Value = '325.51'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("D7").Value = '0.4826'
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").Value = '0.4078'
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("D9").Value = '0.08157'
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("D10").Value = '1.013'
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("E11").Value = '  +4.68%  '
$ws.Range("D12").Value = '1.901.22'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = '6.021'
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("D14").Value = '7.113'
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").Value = '90.43'
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").Value = '0.06788'
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").Value = '17.73'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").Value = '1.006'
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("D21").Value = '29.564.54'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = '5.621'
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("D23").Value = '11.80'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = '2.173'
$ws.Range("E24").Value = '  -1.42%  '
$ws.Range("D25").Value = '2.161.64'
$ws.Range("E25").Value = '  +1.19%  '
$ws.Range("D26").Value = '154.78'
$ws.Range("E26").Value = '  +0.74%  '
$ws.Range("E27").Value = '  +1.20%  '
$ws.Range("D28").Value = '6.291'
$ws.Range("E28").Value = '  +8.44%  '
$ws.Range("D29").Value = '2.102'
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").Value = '119.80'
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").Value = '1.027'
$ws.Range("E31").Value = '  -2.84%  '
$ws.Range("D32").Value = '0.09562'
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").Value = '5.534'
$ws.Range("E33").Value = '  +2.63%  '
$ws.Range("D34").Value = '3.564'
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("D35").Value = '1.393'
$ws.Range("E35").Value = '  -2.18%  '
$ws.Range("D36").Value = '0.02269'
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("D37").Value = '0.06111'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = '1.175'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").Value = '10.85'
$ws.Range("E39").Value = '  +6.92%  '
$ws.Range("D40").Value = '0.5949'
$ws.Range("E40").Value = '  +1.03%  '
$ws.Range("D41").Value = '7.944'
$ws.Range("E41").Value = '  -4.78%  '
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("D43").Value = '2.457'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '1.283'
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").Value = '0.07731'
$ws.Range("E45").Value = '  -3.32%  '
$ws.Range("D46").Value = '12.38'
$ws.Range("E46").Value = '  +1.66%  '
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").Value = '1.951'
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("D49").Value = '115.45'
$ws.Range("E49").Value = '  +1.73%  '
$ws.Range("D50").Value = '72.71'
$ws.Range("E50").Value = '  +1.45%  '
$ws.Range("D51").Value = '1.052'
$ws.Range("E51").Value = '  +1.70%  '
